$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B,C,D,E,G values per row (regenerated s_val data filtering save games)
$data = @{
    2  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    3  = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 13.86384647080068, 20.15985084044064)
    4  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    5  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    6  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    7  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    8  = @(3.272327238179451, 9.983522426115931, 3.223369029078222, 13.86384647080068, 30.34306516417429)
    9  = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 13.86384647080068, 20.15985084044064)
    10 = @(3.272327238179451, 9.983522426115931, 0.7210945179870265, 13.86384647080068, 27.84079065308309)
    11 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 17.65757632934944)
    12 = @(0.04172184405617529, 0.002658071450198252, 0.1496068669990043, 0.5333859586016987, 0.7273727411070765)
    13 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    14 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    15 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
